# Sprint4 Backlog Burndown update
# - Update "Desktop Task" rows (4-7): re-estimate hours and assign "Ahmad" as
#   the team member on each, with Week 2 column (E) now tracked as 0.
# - Move the active selection to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - "Implement and create view for displaying assigned tasks and available tasks"
$ws.Range("C4").Value2 = 1
$ws.Range("D4").Value2 = 1
$ws.Range("E4").Value2 = 0
$ws.Range("F4").Value2 = "Ahmad"

# Row 5 - "Implement and create view for displaying selected task details..."
$ws.Range("C5").Value2 = 1
$ws.Range("D5").Value2 = 1
$ws.Range("E5").Value2 = 0
$ws.Range("F5").Value2 = "Ahmad"

# Row 6 - "Implement functionality to update stage for a task"
$ws.Range("C6").Value2 = 0.5
$ws.Range("D6").Value2 = 1
$ws.Range("E6").Value2 = 0
$ws.Range("F6").Value2 = "Ahmad"

# Row 7 - "Implement functionality to self-assign/unassign to/from a task"
$ws.Range("C7").Value2 = 0.5
$ws.Range("D7").Value2 = 1
$ws.Range("E7").Value2 = 0
$ws.Range("F7").Value2 = "Ahmad"

# Move the selection / active cell to C5 (matches the saved view state).
$ws.Range("C5").Select()
